# A new weekly price observation was inserted at row 36 of the data table.
# Inserting a row shifts every existing record (previously rows 36-62) down
# by one position, so the former last row (62) becomes row 63, and the
# sheet's dimension grows from A1:R62 to A1:R63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36, pushing rows 36-62 down to 37-63.
$ws.Rows("36").Insert()

# Populate the newly inserted row 36 with the new record's data.
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44606
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 60
$ws.Cells.Item(36, 11).Value = 24000
$ws.Cells.Item(36, 12).Value = 25000
$ws.Cells.Item(36, 13).Value = 24500
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(36, 16).Value = 980
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
